$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Absent" column (H) values to reflect the consolidated report
$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 0
